$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the highlight (light-blue fill) from C5 by copying the already
# "un-highlighted" (white) format from B3, which uses the same font (bold grey)
# and alignment, only the fill differs.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Remove the highlight (light-blue fill) from E10 by copying the already
# "un-highlighted" (white) format from E11 (same font/no-alignment style).
$ws.Range("E11").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Clear the resolved note text in F10 ("Tiene que darme acceso a su drive")
$ws.Range("F10").ClearContents() | Out-Null

$excel.CutCopyMode = 0

# Update the saved cursor/selection position
$ws.Range("F3:F4").Select() | Out-Null
